$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers for reference:
# A: Sending cluster, B: Ligand symbol, C: Receptor symbol, D: Target cluster
# E..T: metrics

# New table data (rows 2-10), values per diff.
# Cluster labels: ECs, FAPs, sCs (new cluster) ; Ligand/Receptor symbols stay Ptgs2 / Cav1

$data = @(
    @{ A="ECs";  D="ECs";  E=3; F=1; G=523.0380759999999; H=1569.114228;      I=0.9706650893210215;   J=0.9706650893210214;   K=3; L=1; M=409.6166503333334;  N=1228.849951;     O=0.6234125531262766;    P=0.6234125531262766;    Q=214245.1046879114;   R=1928205.942191203;  S=0.6051248015641634;    T=0.6051248015641633 },
    @{ A="ECs";  D="FAPs"; E=3; F=1; G=523.0380759999999; H=1569.114228;      I=0.9706650893210215;   J=0.9706650893210214;   K=3; L=1; M=56.495384;          N=169.486152;      O=0.08598266586728959;   P=0.08598266586728959;   Q=29549.23695024118;   R=265943.1325521707;  S=0.08346037204413219;   T=0.08346037204413218 },
    @{ A="ECs";  D="sCs";  E=3; F=1; G=523.0380759999999; H=1569.114228;      I=0.9706650893210215;   J=0.9706650893210214;   K=3; L=1; M=190.9434713333333;  N=572.830414;      O=0.2906047810064339;    P=0.2906047810064338;    Q=99870.70587094782;   R=898836.3528385303;  S=0.2820799157127261;    T=0.2820799157127259 },
    @{ A="FAPs"; D="ECs";  E=3; F=1; G=15.69136466666667; H=47.074094;        I=0.02912036538949551;  J=0.02912036538949551;  K=3; L=1; M=409.6166503333334;  N=1228.849951;     O=0.6234125531262766;    P=0.6234125531262766;    Q=6427.444233918823;   R=57846.99810526941;  S=0.01815400133543546;   T=0.01815400133543546 },
    @{ A="FAPs"; D="FAPs"; E=3; F=1; G=15.69136466666667; H=47.074094;        I=0.02912036538949551;  J=0.02912036538949551;  K=3; L=1; M=56.495384;          N=169.486152;      O=0.08598266586728959;   P=0.08598266586728959;   Q=886.4896723273654;   R=7978.407050946289;  S=0.002503846647218377;  T=0.002503846647218376 },
    @{ A="FAPs"; D="sCs";  E=3; F=1; G=15.69136466666667; H=47.074094;        I=0.02912036538949551;  J=0.02912036538949551;  K=3; L=1; M=190.9434713333333;  N=572.830414;      O=0.2906047810064339;    P=0.2906047810064338;    Q=2996.163639410547;   R=26965.47275469492;  S=0.008462517406841679;  T=0.008462517406841678 },
    @{ A="sCs";  D="ECs";  E=1; F=0.3333333333333333; G=0.1156066666666667; H=0.34682; I=0.0002145452894831037; J=0.0002145452894831036; K=3; L=1; M=409.6166503333334;  N=1228.849951;     O=0.6234125531262766;    P=0.6234125531262766;    Q=47.35441555620223;   R=426.1897400058201;  S=0.0001337502266778777; T=0.0001337502266778777 },
    @{ A="sCs";  D="FAPs"; E=1; F=0.3333333333333333; G=0.1156066666666667; H=0.34682; I=0.0002145452894831037; J=0.0002145452894831036; K=3; L=1; M=56.495384;          N=169.486152;      O=0.08598266586728959;   P=0.08598266586728959;   Q=6.531243026293334;   R=58.78118723664;     S=0.00001844717593902662; T=0.00001844717593902662 },
    @{ A="sCs";  D="sCs";  E=1; F=0.3333333333333333; G=0.1156066666666667; H=0.34682; I=0.0002145452894831037; J=0.0002145452894831036; K=3; L=1; M=190.9434713333333;  N=572.830414;      O=0.2906047810064339;    P=0.2906047810064338;    Q=22.07433824260889;   R=198.66904418348;    S=0.0000623478868661993;  T=0.00006234788686619929 }
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = "Ptgs2"
    $ws.Cells.Item($row, 3).Value = "Cav1"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
    $row = $row + 1
}
